# Update odds data for Jogos_da_Semana_FlashScore_2024-11-28.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("J2").Value = 2.88
$ws.Range("S2").Value = 1.47
$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.8
$ws.Range("AG2").Value = 351
$ws.Range("AH2").Value = 9
$ws.Range("AY2").Value = 29

# --- Row 3 updates ---
$ws.Range("G3").Value = 1.45
$ws.Range("H3").Value = 4.5
$ws.Range("I3").Value = 6.5
$ws.Range("J3").Value = 2
$ws.Range("L3").Value = 6.5
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73
$ws.Range("AD3").Value = 8.5
$ws.Range("AH3").Value = 15
$ws.Range("AJ3").Value = 19
$ws.Range("AK3").Value = 67
$ws.Range("AN3").Value = 3.4
$ws.Range("AU3").Value = 9
$ws.Range("AW3").Value = 8
$ws.Range("AX3").Value = 34
$ws.Range("BB3").Value = 301

# --- Row 7 update ---
$ws.Range("N7").Value = 9

# --- Remove the two rows (old rows 8 and 9 - "Al Feiha vs Al Orubah" and
#     "Al Okhdood vs Al Kholood") so that the former row 10
#     ("Al Fateh vs Al Riyadh") shifts up to become the new row 8. ---
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Rows.Item(8).EntireRow.Delete()
